{"js": "// Update the \"Version Control\" table in the Use Case Description document.\n// Three cells in the history table change who is listed as responsible\n// (\u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a) / reviewer (\u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08) for a given version row, while keeping\n// the existing run-level formatting (font/size/etc.) untouched by editing\n// the text of each existing run in place rather than re-creating runs.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each entry: [rowIndex, columnIndex, oldFirstPart, newFirstPart, oldSecondPart, newSecondPart]\n// rowIndex/columnIndex use the logical (merged-cell-aware) indices that\n// Word exposes, matching the table's own header row:\n//   col 3 = \"\u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a\" (responsible), col 4 = \"\u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08\" (reviewer)\nconst edits = [\n  // version 2.8.1 row -> \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08: \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\n  { row: 5, col: 4, oldA: \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", newA: \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\", oldB: \"(SP)\", newB: \" (TL)\" },\n  // version 1.4.1 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM) -> \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\n  { row: 6, col: 3, oldA: \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\", newA: \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", oldB: \" (DM)\", newB: \"(SP)\" },\n  // version 1.4.1 row -> \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08: \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\n  { row: 6, col: 4, oldA: \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", newA: \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\", oldB: \"(SP)\", newB: \" (TL)\" },\n  // version 1.2.2 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\n  { row: 7, col: 3, oldA: \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\", newA: \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\", oldB: \" (DM)\", newB: \" (TL)\" },\n];\n\nfor (const e of edits) {\n  const cell = table.getCell(e.row, e.col);\n\n  const firstResults = cell.body.search(e.oldA, { matchCase: true, matchWholeWord: false });\n  firstResults.load(\"items\");\n  await context.sync();\n  if (firstResults.items.length > 0) {\n    firstResults.items[0].insertText(e.newA, Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  const secondResults = cell.body.search(e.oldB, { matchCase: true, matchWholeWord: false });\n  secondResults.load(\"items\");\n  await context.sync();\n  if (secondResults.items.length > 0) {\n    secondResults.items[0].insertText(e.newB, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Update the \"Version Control\" table in the Use Case Description document.\n# Three version-history rows change who is listed as responsible\n# (\u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a) / reviewer (\u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08), while the existing run-level\n# formatting (font/size/etc.) is preserved because we only overwrite the\n# text of precise sub-ranges located inside the existing runs, rather than\n# replacing whole paragraphs/cells.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Replace-InCell($table, $row, $col, $oldText, $newText) {\n    $cell = $table.Cell($row, $col)\n    $rng = $cell.Range\n    $full = $rng.Text\n    $idx = $full.IndexOf($oldText)\n    if ($idx -ge 0) {\n        $s = $rng.Start + $idx\n        $e = $s + $oldText.Length\n        $target = $d.Range($s, $e)\n        $target.Text = $newText\n    }\n}\n\n# COM table indices are 1-based and count every physical grid cell\n# (rows 1-4 are merged header rows with 2 cells each; the version-history\n# rows start at row 6 and use 5 physical columns: \u0e40\u0e27\u0e2d\u0e23\u0e4c\u0e0a\u0e31\u0e19, \u0e27\u0e31\u0e19\u0e17\u0e35\u0e48,\n# \u0e23\u0e32\u0e22\u0e25\u0e30\u0e40\u0e2d\u0e35\u0e22\u0e14, \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a (col 4), \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08 (col 5)).\n\n# Version 2.8.1 row -> \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08: \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\nReplace-InCell $t 6 5 \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\nReplace-InCell $t 6 5 \"(SP)\" \" (TL)\"\n\n# Version 1.4.1 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM) -> \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\nReplace-InCell $t 7 4 \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\" \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \"\nReplace-InCell $t 7 4 \" (DM)\" \"(SP)\"\n\n# Version 1.4.1 row -> \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08: \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\nReplace-InCell $t 7 5 \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\nReplace-InCell $t 7 5 \"(SP)\" \" (TL)\"\n\n# Version 1.2.2 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\nReplace-InCell $t 8 4 \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\nReplace-InCell $t 8 4 \" (DM)\" \" (TL)\"\n"}
